# B6-PowerPoint.pptx edit script
#
# 1) Three tables (on the slides holding the "Table_0" styled
#    comparison tables) get their table style switched from
#    {AE9A0C15-DED0-4EBE-8D75-D300FE92FE6F} to
#    {D2B2C319-B9BA-4F30-A170-93EF062D308B}.
# 2) The presentation's theme colour scheme is changed from the
#    "Red Violet" / Integral palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newStyleId = "{D2B2C319-B9BA-4F30-A170-93EF062D308B}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $true)
        }
    }
}

# --- 2. Swap the theme colour scheme back to the default Office palette ---
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i]
}
